$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "balance" (H4:K4 merged) and "transactions" (N4) text values
$ws.Range("H4").Value = "0:1"
$ws.Range("N4").Value = "0:2"

# Update the "price" (L4:M4 merged) and the total (K5:N5 merged) numeric values
$ws.Range("L4").Value = 19
$ws.Range("K5").Value = 19
